# Update the "想去人数" (wanted-to-go count) figures for two events.
# These values live in both the "展览" sheet and the "全部类型" sheet,
# which mirror the same data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 196
    $ws.Range("F3").Value = 102
}
